$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detectors")

# Delete row 14 ("801 PS ") and shift the remaining rows up.
$ws.Rows.Item(14).Delete()

# Update the active selection to match the post-edit state.
$ws.Activate()
$ws.Range("A14").Select()
